# Updated cryptos list - refresh Price/Volume(1h) columns, and for a few
# rows the Coin/Link identity moved to a different rank row.
# Cells whose new value still "looks like" a plain number (e.g. "303.72")
# are forced to Text format first so Excel keeps them as strings (matching
# the original inlineStr cells) instead of silently converting them to
# floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.878.59'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '2.271.98'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.72'
$ws.Range('E5').Value = '  +0.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.15'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('E7').Value = '  +1.57%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.65'
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.66'
$ws.Range('E11').Value = '  -1.58%  '
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('E13').Value = '  -1.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.69'
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('D15').Value = '2.624.48'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.32'
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('D17').Value = '2.276.97'
$ws.Range('E17').Value = '  +0.80%  '
$ws.Range('E18').Value = '  +2.99%  '
$ws.Range('D19').Value = '41.793.81'
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.88'
$ws.Range('E20').Value = '  +2.96%  '
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.94'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.29'
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '244.59'
$ws.Range('E24').Value = '  +1.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.58'
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('E26').Value = '  +3.01%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('E28').Value = '  +0.96%  '
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('E30').Value = '  -5.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.14'
$ws.Range('E31').Value = '  +2.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '160.65'
$ws.Range('E32').Value = '  +1.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.27'
$ws.Range('E33').Value = '  +1.16%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0743'
$ws.Range('E35').Value = '  +0.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.03'
$ws.Range('E36').Value = '  -1.49%  '
$ws.Range('E37').Value = '  +1.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.90'
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('E40').Value = '  +0.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.80'
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.94'
$ws.Range('E42').Value = '  -1.55%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.017.54'
$ws.Range('E43').Value = '  -2.23%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.71'
$ws.Range('E44').Value = '  -3.52%  '
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.47'
$ws.Range('E46').Value = '  +3.55%  '
$ws.Range('E47').Value = '  +7.77%  '
$ws.Range('E48').Value = '  -2.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.56'
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.17'
$ws.Range('E50').Value = '  +3.46%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.51'
$ws.Range('E51').Value = '  -1.04%  '
